$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-modified date serial number for each
# data row (rows 2 through 401). Bump it from 45202 (2023-10-03) to
# 45203 (2023-10-04) for every row, matching the upstream refresh.
$ws.Range("C2:C401").Value = 45203
